$wb = $excel.ActiveWorkbook

# --- open_tasks: remove the "refactor ffn, dc, m_nb" / Robin row ---
$ws1 = $wb.Worksheets.Item("open_tasks")
$ws1.Rows.Item(9).Delete()

# --- done_tasks: update existing numbers and add new rows ---
$ws2 = $wb.Worksheets.Item("done_tasks")

# Row 8 (Preference Extraction) updates
$ws2.Range("B8").Value = 210
$ws2.Range("F8").Value = 180

# Row 9 (Design of FSM) updates
$ws2.Range("C9").Value = 90
$ws2.Range("D9").Value = 210
$ws2.Range("F9").Value = 30

# Insert new row 10: Implementing FSM
$ws2.Rows.Item(9).Copy()
$ws2.Rows.Item(10).Insert()
$ws2.Range("A10").Value = "Implementing FSM"
$ws2.Range("B10").Value = ""
$ws2.Range("C10").Value = ""
$ws2.Range("D10").Value = ""
$ws2.Range("E10").Value = 90
$ws2.Range("F10").Value = 30

# Insert new row 11: refactor ffn, dc, m_nb
$ws2.Rows.Item(11).Copy()
$ws2.Rows.Item(11).Insert()
$ws2.Range("A11").Value = "refactor ffn, dc, m_nb"
$ws2.Range("B11").Value = ""
$ws2.Range("C11").Value = ""
$ws2.Range("D11").Value = ""
$ws2.Range("E11").Value = ""
$ws2.Range("F11").Value = 30

# Insert new row 12: empty row
$ws2.Rows.Item(12).Copy()
$ws2.Rows.Item(12).Insert()
$ws2.Range("A12").Value = ""
$ws2.Range("B12").Value = ""
$ws2.Range("C12").Value = ""
$ws2.Range("D12").Value = ""
$ws2.Range("E12").Value = ""
$ws2.Range("F12").Value = ""
